$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing quantity for the first BOM line (1000 units -> 500 units) ---
$ws.Range("G4").Value = "500 units"

# --- Add a new BOM line (row 6): COC_CASE_CHST_ONE_HOLD ---
$ws.Range("C6").Value = "COC_CASE_CHST_ONE_HOLD"
$ws.Range("D6").Value = "Adhesive Circles`n4 per case"
$ws.Range("F6").Value = "Savanna Warehouse"
$ws.Range("G6").Value = "3000 units"

# Match the formatting used by the other populated rows (center-aligned, wrapped where needed)
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").WrapText = $true
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("G6").HorizontalAlignment = -4108

# F5 (Supplier on the 2nd BOM line) becomes center-aligned like the rest of the column
$ws.Range("F5").HorizontalAlignment = -4108

# --- Row heights for the header rows shrink slightly ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 15

# --- Column width adjustments (A & B become one uniform width, others re-sized) ---
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 13.833333333333334
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(7).ColumnWidth = 19.833333333333332

# --- Selection moves to the merged A2:B2 header cell ---
$ws.Range("A2:B2").Select() | Out-Null
